$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column B (Total) updates
$ws.Range("B3").Value = 8497.212576942569
$ws.Range("B4").Value = 3548.124252375236
$ws.Range("B5").Value = 2665.504189041093
$ws.Range("B6").Value = 5601.561525342468
$ws.Range("B7").Value = 7381.684523287677
$ws.Range("B8").Value = 11460.53568630137
$ws.Range("B9").Value = 16057.41333150689

# Column D (Community) updates
$ws.Range("D3").Value = 576.0166857515222
$ws.Range("D4").Value = 427.3954172529718
$ws.Range("D5").Value = 32.60295890410961
$ws.Range("D6").Value = 149.0417726027392
$ws.Range("D7").Value = 523.5343575342463
$ws.Range("D8").Value = 1040
$ws.Range("D9").Value = 1040.00004109589

# Totals / ratios updates
$ws.Range("F10").Value = 12235986.11079729
$ws.Range("G11").Value = 0.7143008952600087
$ws.Range("F12").Value = 829464.0274821916
$ws.Range("G12").Value = 0.06778889906962672
$ws.Range("G13").Value = 0.2179102056703646
